$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad" date) for rows 2-7 from 2023-09-16 (45185) to 2023-10-05 (45204)
$newDate = Get-Date -Year 2023 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
